$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp note in A1
$ws.Range("A1").Value = "Datos actualizados a 30 de Marzo de 2020 a las 04:50"

# Row 4 (A4 = countries index 8)
$ws.Range("B4").Value = 142070
$ws.Range("C4").Value = 23
$ws.Range("E4").Value = 135027

# Row 6 (A6 = countries index 10)
$ws.Range("B6").Value = 81470
$ws.Range("C6").Value = 31
$ws.Range("D6").Value = 75700
$ws.Range("E6").Value = 2466
$ws.Range("F6").Value = 633
$ws.Range("G6").Value = 4
$ws.Range("H6").Value = 3304

# Row 15 (A15 = countries index 19)
$ws.Range("B15").Value = 9661
$ws.Range("C15").Value = 78
$ws.Range("D15").Value = 5228
$ws.Range("E15").Value = 4275
$ws.Range("G15").Value = 6
$ws.Range("H15").Value = 158

# Row 63 (A63 = countries index 67)
$ws.Range("D63").Value = 63
$ws.Range("E63").Value = 450
$ws.Range("F63").Value = 2
